$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 1, shifting all existing data
# (header + 60 data rows) down by one row.
$ws.Rows.Item(1).Insert()
